$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: the old "_GoBack" bookmark currently sits between the runs
# "B6" and "5_S3_LisezMoi_" near the end of the document. Remove it and
# merge those two runs' text into a single run "B65_S3_LisezMoi_"
# (they already share identical formatting, so replacing the text of the
# first run and deleting the now-redundant text of the second run
# achieves a clean merge while keeping the shared rPr formatting).
# -----------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
}

$mergeRange = $d.Content
$foundB6 = $mergeRange.Find.Execute("B6", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
if ($foundB6) {
    # Rewrite the first run's text in place (keeps its run formatting).
    $mergeRange.Text = "B65_S3_LisezMoi_"
    $tailStart = $mergeRange.End

    # The old second run's text ("5_S3_LisezMoi_") immediately follows;
    # remove it now that its text has been folded into the first run.
    $tailRange = $d.Range($tailStart, $d.Content.End)
    $foundTail = $tailRange.Find.Execute("5_S3_LisezMoi_", $true, $false, $false, `
                                          $false, $false, $true, 1, $false, "", 0)
    if ($foundTail) {
        $tailRange.Delete()
    }
}

# -----------------------------------------------------------------------
# Part 2: add a new "_GoBack" bookmark at the very end of the paragraph
# that ends with "... l'utilisation de l'application. " (right after the
# last run, before the paragraph mark).
#
# Note: this runtime mis-places a zero-length bookmark range that lands
# exactly on a paragraph-mark position. Work around that by temporarily
# inserting a one-character placeholder right at the target position
# (which shifts the paragraph mark forward so the target position is no
# longer "on" it), adding the bookmark there, then deleting the
# placeholder (it sits after the now-placed bookmark, so removing it
# does not move the bookmark).
# -----------------------------------------------------------------------

$paraRange = $d.Content
$foundPara = $paraRange.Find.Execute("nécessaire à l’utilisation de l’application. ", `
                                      $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if ($foundPara) {
    $targetPos = $paraRange.End

    $placeholder = $d.Range($targetPos, $targetPos)
    $placeholder.InsertAfter("Z")

    $bookmarkRange = $d.Range($targetPos, $targetPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    $placeholderRange = $d.Range($targetPos, $targetPos + 1)
    $placeholderRange.Delete()
}
